# Add git diff usage
# Inserts a new row (row 14) on Sheet1 describing how to diff a file between
# two revisions, right after the existing "diff" row, and updates the sheet
# font from the CJK default (新細明體) to Calibri, matching the authored
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row after row 13 ("Git" / "diff" / staged-files tip), pushing
# everything below down by one.
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = "Git"
$ws.Cells.Item(14, 2).Value = "diff of a file between commits"
$ws.Cells.Item(14, 3).Value = "git diff <revision_1>:<file_1> <revision_2>:<file_2>`nfor example: git diff master:pom.xml d44ac61:pom.xml"

# Match the formatting (style 5 = wrap-text body cells) used by the rest of
# the table's data rows.
$ws.Cells.Item(14, 1).Style = $ws.Cells.Item(13, 1).Style
$ws.Cells.Item(14, 2).Style = $ws.Cells.Item(13, 2).Style
$ws.Cells.Item(14, 3).Style = $ws.Cells.Item(13, 3).Style

$ws.Rows.Item(14).RowHeight = 26.25

# The workbook-wide font moved from the CJK "新細明體" default to Calibri.
$wb.Worksheets.Item("Sheet1").Cells.Font.Name = "Calibri"
$wb.Worksheets.Item("Sheet2").Cells.Font.Name = "Calibri"
$wb.Worksheets.Item("Sheet3").Cells.Font.Name = "Calibri"
